# Update "想去人数" (want-to-go count) values in column F across the
# workbook's four sheets, matching the new scrape output committed as
# "Update gh-pages to output generated at 456a3b4".
#
# Sheet order (per workbook.xml): 1=展览 (Exhibitions), 2=演出 (Shows),
# 3=本地生活 (Local life), 4=全部类型 (All types / combined view).

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibitions)
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 4568
$ws.Range("F4").Value = 447
$ws.Range("F5").Value = 3683
$ws.Range("F6").Value = 1062
$ws.Range("F9").Value = 371
$ws.Range("F10").Value = 366
$ws.Range("F11").Value = 2550
$ws.Range("F12").Value = 1285
$ws.Range("F13").Value = 39
$ws.Range("F14").Value = 1968
$ws.Range("F16").Value = 20
$ws.Range("F19").Value = 64
$ws.Range("F20").Value = 10540
$ws.Range("F21").Value = 6115
$ws.Range("F25").Value = 218
$ws.Range("F30").Value = 182
$ws.Range("F31").Value = 862
$ws.Range("F32").Value = 3566
$ws.Range("F36").Value = 128
$ws.Range("F37").Value = 275
$ws.Range("F39").Value = 253
$ws.Range("F40").Value = 4867
$ws.Range("F42").Value = 1149
$ws.Range("F43").Value = 168
$ws.Range("F44").Value = 193
$ws.Range("F45").Value = 108

# Sheet 2: 演出 (Shows)
$ws = $wb.Worksheets.Item(2)
$ws.Range("F15").Value = 3602

# Sheet 3: 本地生活 (Local life)
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 8841
$ws.Range("F3").Value = 449
$ws.Range("F4").Value = 1661

# Sheet 4: 全部类型 (All types / combined view)
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 449
$ws.Range("F4").Value = 4568
$ws.Range("F6").Value = 447
$ws.Range("F7").Value = 3683
$ws.Range("F10").Value = 2550
$ws.Range("F14").Value = 1285
$ws.Range("F15").Value = 39
$ws.Range("F16").Value = 20
$ws.Range("F20").Value = 64
$ws.Range("F21").Value = 10540
$ws.Range("F22").Value = 3602
$ws.Range("F27").Value = 218
$ws.Range("F32").Value = 182
$ws.Range("F33").Value = 862
$ws.Range("F34").Value = 3566
$ws.Range("F36").Value = 128
$ws.Range("F37").Value = 275
$ws.Range("F40").Value = 253
$ws.Range("F41").Value = 4867
$ws.Range("F43").Value = 1149
$ws.Range("F44").Value = 168
